$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = 53
$ws.Range("G1").Value = 1.8
